$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "300.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.70%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "38.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "8.59%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.984"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-3.39%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.54%"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-5.94%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.965"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.88%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.994"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.28%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9164"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.61%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09091"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-8.63%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1795"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.01%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08445"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.70%"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "6.75%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09939"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.16%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001493"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.48%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005681"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.65%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.476"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.36%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.99%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2.87%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1317"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.21%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.567"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "6.30%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2235"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.94%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04656"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.64%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.23%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004439"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.40%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.01%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004757"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "40.11%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01737"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-3.17%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04688"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.24%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007895"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.39%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1388"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.73%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007687"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "6.68%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002302"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "11.08%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009777"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.41%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006030"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-1.35%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.03%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.726"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "189.28%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "34.83%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.03%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.03%"
